# Integrating ExtentReport for parallel tests
# Trim the "TestData" sheet's sample rows down to a smaller, faster data
# set (fewer AddCustomerTest / OpenAccountTest rows) and update the
# selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Wipe the previous (larger) table so stale cells don't linger past the
# new, smaller extent (old data ran through row 20).
$ws.Range("A1:E20").ClearContents()

# --- AddCustomerTest block ---------------------------------------------
$ws.Range("A1").Value = "AddCustomerTest"

$ws.Range("A2").Value = "Runmode"
$ws.Range("B2").Value = "firstname"
$ws.Range("C2").Value = "lastname"
$ws.Range("D2").Value = "postcode"
$ws.Range("E2").Value = "browser"

$ws.Range("A3").Value = "Y"
$ws.Range("B3").Value = "Ivan"
$ws.Range("C3").Value = "Ivanov"
$ws.Range("D3").Value = "e3r4t5"
$ws.Range("E3").Value = "chrome"

$ws.Range("A4").Value = "Y"
$ws.Range("B4").Value = "Petr"
$ws.Range("C4").Value = "Petrov"
$ws.Range("D4").Value = "2af4g5"
$ws.Range("E4").Value = "firefox"

# --- OpenAccountTest block ----------------------------------------------
$ws.Range("A6").Value = "OpenAccountTest"

$ws.Range("A7").Value = "Runmode"
$ws.Range("B7").Value = "customer"
$ws.Range("C7").Value = "currency"
$ws.Range("D7").Value = "browser"

$ws.Range("A8").Value = "Y"
$ws.Range("B8").Value = "Ivan Ivanov"
$ws.Range("C8").Value = "Dollar"
$ws.Range("D8").Value = "chrome"

$ws.Range("A9").Value = "Y"
$ws.Range("B9").Value = "Petr Petrov"
$ws.Range("C9").Value = "Rupee"
$ws.Range("D9").Value = "firefox"

# Update the selection to match the post-edit workbook state.
$ws.Activate() | Out-Null
$ws.Range("F12").Select() | Out-Null
